# Updates the "NEW" mapa_interactivo dataset to the latest export:
#  - inserts a new incident report (Caso 6036) that sorts in by date between
#    the existing 1/21/2025 and 3/7/2025 rows, shifting rows 7-12 down to 8-13
#    (the previously-last row in that date run, Caso 5458, drops out of the
#    visible window as a result)
#  - removes the resolved Caso 5850 (BLANCO ENCALADA) entry, shifting rows
#    34-43 up to 33-42, and appends the newest report (Caso 6020) at row 43

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($ws, $rowNum, $vals, $mx, $my)
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($rowNum, $col)
        # Columns A-L are always stored as text in this sheet (even the
        # numeric-looking Caso/OT/Comuna values and the dates in column B),
        # so force text formatting before writing to stop Excel from
        # auto-coercing them into numbers / date serials.
        $cell.NumberFormat = "@"
        $cell.Value2 = $vals[$i]
    }
    # Coordenada_X / Coordenada_Y (M/N) are genuine numbers.
    $mCell = $ws.Cells.Item($rowNum, 13)
    $mCell.NumberFormat = "General"
    $mCell.Value2 = $mx
    $nCell = $ws.Cells.Item($rowNum, 14)
    $nCell.NumberFormat = "General"
    $nCell.Value2 = $my
}

Set-RowData $ws 7 @("6036", "2/24/2025", "MEDRANO 1715", "14", "803608181", "NEW", "Pendiente", "", "0", "Cambio", "Nodo Teco", "Pasante") -58.418236 -34.589859
Set-RowData $ws 8 @("5037", "3/7/2025", "Monroe 3605", "12", "803825082", "NEW", "Pendiente", "Columna inclinada", "0", "Aplomo", "Sin equipos", "Pasante") -58.471774 -34.565411
Set-RowData $ws 9 @("5053", "3/11/2025", "MENDOZA 1153", "13", "803969314", "NEW", "Pendiente", "Poste", "0", "Cambio", "Sin equipos", "Poste") -58.44463 -34.553354
Set-RowData $ws 10 @("803969337", "3/11/2025", "Rousseau 2087", "15", "803969337", "NEW", "Pendiente", " poste de Telefonia por caer", "0", "Cambio", "Sin equipos", "Poste") -58.47678 -34.601336
Set-RowData $ws 11 @("2125", "3/27/2025", "DIAZ COLODRERO 2321", "12", "804309658", "NEW", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante") -58.485065 -34.574269
Set-RowData $ws 12 @("3348", "3/27/2025", "ESTOMBA 2626", "12", "804309704", "NEW", "Pendiente", "Picada", "1", "Desmonte", "Sin equipos", "Pasante") -58.47538 -34.566015
Set-RowData $ws 13 @("3430", "4/1/2025", "MONROE 3838", "12", "804468442", "NEW", "Pendiente", "Reparar rienda", "1", "Tensor", "Sin equipos", "Terminal") -58.473659 -34.566979
Set-RowData $ws 33 @("5855", "5/5/2025", "IBERA 4960", "12", "805676619", "NEW", "Pendiente", "Cambiar columna corroída en la base tiene 2 CDO ", "1", "Cambio", "Sin equipos", "Pasante") -58.489018 -34.566163
Set-RowData $ws 34 @("5847", "5/8/2025", "JURAMENTO 5321", "12", "805791839", "NEW", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante") -58.485193 -34.579621
Set-RowData $ws 35 @("5802", "5/13/2025", "MACHAIN 4516", "12", "806926363", "NEW", "Pendiente", "Poste podrido", "1", "Cambio", "Sin equipos", "Poste") -58.49243 -34.551559
Set-RowData $ws 36 @("5826", "5/19/2025", "ALBARELLOS AV. 3180", "12", "806926466", "NEW", "Pendiente", "Columna (metal) inclinada", "1", "Aplomo", "Sin equipos", "Terminal") -58.513552 -34.579829
Set-RowData $ws 37 @("5825", "5/19/2025", "PAZ, GRAL. AV. 5602", "12", "806926472", "NEW", "Pendiente", "reparar rienda cortada - ver foto no es la misma que albarellos", "1", "Tensor", "Sin equipos", "Terminal") -58.513685 -34.579838
Set-RowData $ws 38 @("806926510", "5/22/2025", "Paz Soldan 4991", "15", "806926510", "NEW", "Pendiente", "Poste inclinado, cambiar o desmontar", "1", "Cambio", "Sin equipos", "Poste") -58.468466 -34.594154
Set-RowData $ws 39 @("5863", "5/27/2025", "QUINTANA 4631", "12", "806975681", "NEW", "Pendiente", "Poste quebrado", "1", "Cambio", "Sin equipos", "Poste") -58.480224 -34.544229
Set-RowData $ws 40 @("5875", "5/27/2025", "MONROE 4096", "12", "806975680", "NEW", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante") -58.476106 -34.568373
Set-RowData $ws 41 @("-447", "5/28/2025", "Avenida Balbin 3883", "12", "806975696", "NEW", "Pendiente", "Columna corroida esta desprendida", "1", "Cambio", "Sin equipos", "Pasante") -58.484982 -34.554653
Set-RowData $ws 42 @("807044192", "5/29/2025", "O'Higgins 4379", "13", "807044192", "NEW", "Pendiente", "Poste ", "1", "Desmonte", "Sin equipos", "Poste") -58.468425 -34.54124
Set-RowData $ws 43 @("6020", "6/5/2025", "RAVIGNANI, EMILIO, DR. 2036", "14", "807215465", "NEW", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante") -58.436298 -34.578972
